# Fruta / hortaliza, semanal
# Insert two new weekly records (row 636 and 637) into the "Naranja" price
# sheet, pushing the previously-existing rows 636-710 down to 638-712.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 636 (inserting twice at the same index
# pushes everything that was at 636 downward, leaving two fresh rows).
$ws.Rows.Item(636).Insert()
$ws.Rows.Item(636).Insert()

# New record: Fukumoto / Primera
$row636 = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45124, 16, "Fruta", 100102, "Cítricos", 100102005, "Naranja", "Fukumoto", "Primera", 80, 10000, 10000, 10000, "`$/bandeja 15 kilos granel", "Región de O'Higgins", 667, 15)

# New record: Fukumoto / Segunda
$row637 = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45124, 16, "Fruta", 100102, "Cítricos", 100102005, "Naranja", "Fukumoto", "Segunda", 80, 8000, 8000, 8000, "`$/bandeja 15 kilos granel", "Región de O'Higgins", 533, 15)

for ($i = 0; $i -lt $row636.Count; $i++) {
    $ws.Cells.Item(636, $i + 1).Value = $row636[$i]
}

for ($i = 0; $i -lt $row637.Count; $i++) {
    $ws.Cells.Item(637, $i + 1).Value = $row637[$i]
}
